$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 133
$ws.Range("F8").Value = 594
$ws.Range("F9").Value = 4
$ws.Range("F10").Value = 5
$ws.Range("F12").Value = 10482
$ws.Range("F14").Value = 75
$ws.Range("F16").Value = 1998
$ws.Range("F17").Value = 881
$ws.Range("F20").Value = 189
$ws.Range("F22").Value = 227
$ws.Range("F23").Value = 1137
$ws.Range("F25").Value = 156
$ws.Range("F26").Value = 656
$ws.Range("F27").Value = 63
$ws.Range("F28").Value = 197
$ws.Range("F29").Value = 648
$ws.Range("F30").Value = 2975
$ws.Range("F31").Value = 984
$ws.Range("F32").Value = 712
$ws.Range("F34").Value = 19
$ws.Range("F36").Value = 897
$ws.Range("F37").Value = 10
$ws.Range("F41").Value = 1153
$ws.Range("F43").Value = 87
$ws.Range("F44").Value = 120
$ws.Range("F45").Value = 205
$ws.Range("F47").Value = 3
$ws.Range("F49").Value = 73

$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 4063
$ws.Range("F14").Value = 226

$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 400

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 400
$ws.Range("F7").Value = 133
$ws.Range("F13").Value = 594
$ws.Range("F15").Value = 10482
$ws.Range("F16").Value = 75
$ws.Range("F18").Value = 1998
$ws.Range("F19").Value = 881
$ws.Range("F22").Value = 1137
$ws.Range("F24").Value = 156
$ws.Range("F25").Value = 4063
$ws.Range("F26").Value = 656
$ws.Range("F27").Value = 63
$ws.Range("F28").Value = 197
$ws.Range("F29").Value = 648
$ws.Range("F30").Value = 2975
$ws.Range("F31").Value = 984
$ws.Range("F34").Value = 712
$ws.Range("F35").Value = 19
$ws.Range("F38").Value = 10
$ws.Range("F41").Value = 1153
$ws.Range("F43").Value = 87
$ws.Range("F44").Value = 120
$ws.Range("F45").Value = 205
